$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column L ("median_abschg") holds a copy of column K ("mean_abschg") for
# every data row (2-89). Per the target commit, the median_abschg column
# should be zeroed out across the whole table.
$lastRow = $ws.Cells(1, 12).End(-4121).Row  # xlDown = -4121
if ($lastRow -lt 89) { $lastRow = 89 }

$ws.Range("L2:L$lastRow").Value = 0
